# Generate Report for Handoff
# Update status/handoff info for the two files (52dd6477-... and a3984ac7-...)
# that are now ready for handoff, across the Overview sheet and the
# per-locale (zh-cn / de-de) sheets. Also widen the "Error Detail" column
# on the locale sheets so the new long error messages are readable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: rows 4 (52dd6477-...) and 5 (a3984ac7-...)
#   zh-cn / de-de status columns -> "Ready for handoff"
#   Latest HO Xliff Generate Date -> 2016-08-22 11:49:35
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-08-22 11:49:35"

$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-22 11:49:35"

# ---------------------------------------------------------------------
# zh-cn sheet: rows 4 (52dd6477-...) and 5 (a3984ac7-...)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("H4").Value = "2016-08-22 11:49:30"
$wsZhCn.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b9c5f576b1b5de8a675a6136806fa680e2df24c/e2e/52dd6477-1001-484b-9d5e-4c74afb68037.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8737a983cccf15e5e6b3625707f4ef0da9e5226d/e2e/52dd6477-1001-484b-9d5e-4c74afb68037.md."

$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("H5").Value = "2016-08-22 11:49:30"
$wsZhCn.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b9c5f576b1b5de8a675a6136806fa680e2df24c/e2e/a3984ac7-bb4a-41af-9e85-793069935bc7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8737a983cccf15e5e6b3625707f4ef0da9e5226d/e2e/a3984ac7-bb4a-41af-9e85-793069935bc7.md."

# widen the Error Detail column (P) so the new long messages are readable
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: rows 4 (52dd6477-...) and 5 (a3984ac7-...)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("H4").Value = "2016-08-22 11:49:35"
$wsDeDe.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b9c5f576b1b5de8a675a6136806fa680e2df24c/e2e/52dd6477-1001-484b-9d5e-4c74afb68037.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8737a983cccf15e5e6b3625707f4ef0da9e5226d/e2e/52dd6477-1001-484b-9d5e-4c74afb68037.md."

$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("H5").Value = "2016-08-22 11:49:35"
$wsDeDe.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b9c5f576b1b5de8a675a6136806fa680e2df24c/e2e/a3984ac7-bb4a-41af-9e85-793069935bc7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8737a983cccf15e5e6b3625707f4ef0da9e5226d/e2e/a3984ac7-bb4a-41af-9e85-793069935bc7.md."

# widen the Error Detail column (P) so the new long messages are readable
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
